$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29; this pushes the existing rows 29..112
# down to 30..113 (and the sheet dimension grows to A1:R113).
$ws.Rows("29:29").Insert()

# Populate the newly inserted row 29 with the new weekly record.
$ws.Range("A29").Value = 6
$ws.Range("B29").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C29").Value = "Metropolitana"
$ws.Range("D29").Value = 44453
$ws.Range("D29").NumberFormat = $ws.Range("D30").NumberFormat
$ws.Range("E29").Value = 13
$ws.Range("F29").Value = 100112026
$ws.Range("G29").Value = "Haba"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 800
$ws.Range("K29").Value = 12000
$ws.Range("L29").Value = 13000
$ws.Range("M29").Value = 12562
$ws.Range("N29").Value = "`$/saco 25 kilos"
$ws.Range("O29").Value = "Región de Coquimbo"
$ws.Range("P29").Value = 502
$ws.Range("Q29").Value = 25
$ws.Range("R29").Value = "Hortaliza"
